function Set-Headers($ws) {
    $headers = @("prompt","solution","llm_response","evaluator_response","evaluator_partial_correctness")
    $cols = @("A","B","C","D","E")
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Range($cols[$i] + "1").Value = $headers[$i]
    }
    $full = $ws.Range("A1:E1")
    $full.HorizontalAlignment = -4108
    $full.VerticalAlignment = -4160
    $full.Borders.LineStyle = 1
    $full.Font.Bold = $true
}

$wb = $excel.ActiveWorkbook

# --- Sheet1: o_10 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Set-Headers $ws1

$promptO10 = @'
 Given is the adjacency matrix for a weighted directed graph containing 15 nodes labelled A to O. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node O? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O
 A 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 2 0 3 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0
 F 0 0 0 0 0 0 4 0 0 5 0 0 0 0 0
 G 0 0 0 0 0 4 0 0 0 0 5 0 0 0 0
 H 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0
 I 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0
 J 0 0 0 0 0 2 0 0 5 0 4 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3
 L 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0
 N 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0
    
'@

$solutionO10 = @'
A -> B -> C -> D -> G -> K -> O
'@

$llmRespO10 = @'
The least cost path from node A to node O is A-B-C-D-G-K-O, with a total cost of 10.
'@

$evalRespO10 = @'
invalid input
'@

$evalPartialO10 = @'
7/7
'@

$ws1.Range("A2").Value = $promptO10
$ws1.Range("B2").Value = $solutionO10
$ws1.Range("C2").Value = $llmRespO10
$ws1.Range("D2").Value = $evalRespO10
$ws1.Range("E2").Value = $evalPartialO10

# --- Sheet2: o_20 -------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"
Set-Headers $ws2

$promptO20 = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 2 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 5 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 2 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 5 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 3 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 3
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@

$solutionO20 = @'
A -> F -> K -> P -> U -> V -> W -> X -> Y
'@

$llmRespO20 = @'
Based on the given adjacency matrix for a weighted directed graph, the least cost path from node A to node Y is: 
A -> B -> C -> H -> I -> J -> O -> T -> Y 
Here is the cost of each step in the path:
A -> B : 2 
B -> C : 5 
C -> H : 5 
H -> I : 4 
I -> J : 1 
J -> O : 1 
O -> T : 1 
T -> Y : 3 
So, the total cost for the least path from node A to node Y is 2 + 5 + 5 + 4 + 1 + 1 + 1 + 3 = 22.
'@

$evalRespO20 = @'
invalid input
'@

$evalPartialO20 = @'
0/9
'@

$ws2.Range("A2").Value = $promptO20
$ws2.Range("B2").Value = $solutionO20
$ws2.Range("C2").Value = $llmRespO20
$ws2.Range("D2").Value = $evalRespO20
$ws2.Range("E2").Value = $evalPartialO20

# --- Sheet3: o_20_jumbled ------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"
Set-Headers $ws3

$promptO20J = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 3 0 5 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 1 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 5 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 4 0 1 0 0 0 5 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 4 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 5 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 4
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@

$solutionO20J = @'
A -> F -> G -> H -> I -> N -> O -> T -> Y
'@

$llmRespO20J = @'
Based on the adjacency matrix, we can see that the smallest cost is 1. Following the path with these costs, we find that the least cost path from A to Y is A - F - G - L - M - R - S - T - Y. The total cost of this path is 1 + 2 + 1 + 2 + 1 + 2 + 3 + 1 = 13.
'@

$evalRespO20J = @'
invalid input
'@

$evalPartialO20J = @'
2/9
'@

$ws3.Range("A2").Value = $promptO20J
$ws3.Range("B2").Value = $solutionO20J
$ws3.Range("C2").Value = $llmRespO20J
$ws3.Range("D2").Value = $evalRespO20J
$ws3.Range("E2").Value = $evalPartialO20J

# Restore o_10 as the selected/active sheet (matches original tabSelected="1").
$ws1.Activate()

Write-Host "Workbook updated: 3 sheets (o_10, o_20, o_20_jumbled)"
